$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.107.30"
$ws.Range("E2").Value = "  -0.21%  "
$ws.Range("D3").Value = "1.625.80"
$ws.Range("E3").Value = "  -0.93%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'214.23"
$ws.Range("E5").Value = "  -1.35%  "
$ws.Range("E6").Value = "  +1.76%  "
$ws.Range("E7").Value = "  -0.04%  "
$ws.Range("E8").Value = "  -1.21%  "
$ws.Range("E9").Value = "  +0.22%  "
$ws.Range("D10").Value = "'20.36"
$ws.Range("E10").Value = "  +1.78%  "
$ws.Range("D11").Value = "'0.0848"
$ws.Range("E11").Value = "  +0.04%  "
$ws.Range("D12").Value = "1.629.22"
$ws.Range("E12").Value = "  -1.20%  "
$ws.Range("D13").Value = "'4.14"
$ws.Range("E13").Value = "  -0.20%  "
$ws.Range("E14").Value = "  +0.02%  "
$ws.Range("B15").Value = "Litecoin"
$ws.Range("C15").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D15").Value = "'64.66"
$ws.Range("E15").Value = "  -3.88%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "27.108.38"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "0.0₃0744"
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "'217.22"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("E19").Value = "  -0.02%  "
$ws.Range("D20").Value = "'6.97"
$ws.Range("E20").Value = "  +1.45%  "
$ws.Range("D21").Value = "'4.37"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("E22").Value = "  -6.07%  "
$ws.Range("D23").Value = "'9.07"
$ws.Range("E23").Value = "  -1.53%  "
$ws.Range("D24").Value = "'148.12"
$ws.Range("E24").Value = "  +0.26%  "
$ws.Range("E25").Value = "  +0.03%  "
$ws.Range("D26").Value = "'7.30"
$ws.Range("E26").Value = "  -3.13%  "
$ws.Range("E27").Value = "  -0.53%  "
$ws.Range("D28").Value = "'15.62"
$ws.Range("E28").Value = "  -0.98%  "
$ws.Range("E29").Value = "  -0.27%  "
$ws.Range("E30").Value = "  -1.29%  "
$ws.Range("D31").Value = "'3.37"
$ws.Range("E31").Value = "  -0.26%  "
$ws.Range("E32").Value = "  -1.03%  "
$ws.Range("D33").Value = "1.349.08"
$ws.Range("E33").Value = "  +5.97%  "
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("E35").Value = "  -0.35%  "
$ws.Range("D36").Value = "'0.0178"
$ws.Range("E36").Value = "  +0.15%  "
$ws.Range("E37").Value = "  +1.44%  "
$ws.Range("D38").Value = "'0.857"
$ws.Range("E38").Value = "  +0.31%  "
$ws.Range("E39").Value = "  -0.04%  "
$ws.Range("E40").Value = "  -0.51%  "
$ws.Range("E41").Value = "  +0.21%  "
$ws.Range("D42").Value = "'65.60"
$ws.Range("E42").Value = "  +6.28%  "
$ws.Range("B43").Value = "FraxShare"
$ws.Range("C43").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D43").Value = "'5.24"
$ws.Range("E43").Value = "  -1.31%  "
$ws.Range("B44").Value = "RocketPoolETH"
$ws.Range("C44").Value = "https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth"
$ws.Range("D44").Value = "1.763.71"
$ws.Range("E44").Value = "  -1.08%  "
$ws.Range("D45").Value = "'90.66"
$ws.Range("E45").Value = "  -1.28%  "
$ws.Range("E46").Value = "  +0.89%  "
$ws.Range("D47").Value = "'0.856"
$ws.Range("E47").Value = "  +28.84%  "
$ws.Range("B48").Value = "BabyDogeCoin"
$ws.Range("C48").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D48").Value = "0.0₆0106"
$ws.Range("E48").Value = "  -1.13%  "
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "  +0.05%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").Value = "'0.0995"
$ws.Range("E50").Value = "  +2.18%  "
$ws.Range("B51").Value = "EnergySwap"
$ws.Range("C51").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.59"
$ws.Range("E51").Value = "  -0.38%  "
